# Reverse the order of the comma-separated "Recorded By" entries in column G.
# e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com"
# Single-value cells (no separator) are left unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "") {
        $parts = $val.Split(", ")
        if ($parts.Count -gt 1) {
            $reversed = $parts[($parts.Count - 1)..0]
            $newVal = $reversed -join ", "
            $cell.Value = $newVal
        }
    }
}
